$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row: Date (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P), Precio $/Kg (S)
$data = @{
    2  = @{ D = 44417; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 }
    3  = @{ D = 44435; M = 130; N = 1300; O = 1300; P = 1300; S = 1300 }
    4  = @{ D = 44438; M = 60;  N = 1200; O = 1200; P = 1200; S = 1200 }
    5  = @{ D = 44343; M = 60;  N = 1300; O = 1300; P = 1300; S = 1300 }
    6  = @{ D = 44431; M = 100; N = 1300; O = 1300; P = 1300; S = 1300 }
    7  = @{ D = 44424; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 }
    8  = @{ D = 44476; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 }
    9  = @{ D = 44418; M = 40;  N = 1200; O = 1200; P = 1200; S = 1200 }
    10 = @{ D = 44405; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 }
    11 = @{ D = 44432; M = 30;  N = 1300; O = 1300; P = 1300; S = 1300 }
    12 = @{ D = 44473; M = 120; N = 1200; O = 1200; P = 1200; S = 1200 }
    13 = @{ D = 44357; M = 35;  N = 1000; O = 1000; P = 1000; S = 1000 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $vals.S   # S: Precio $/Kg
}
